# Edit sheet Card24 by admin
# Fill all currently-blank inline-string cells in columns D:L (rows 2-13) with the text "nan",
# leaving cells that already contain a value (checkmarks, dates, numbers, etc.) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Card24")

$cols = @("D","E","F","G","H","I","J","K","L")

for ($r = 2; $r -le 13; $r++) {
    foreach ($col in $cols) {
        $cell = $ws.Range("$col$r")
        $current = $cell.Value()
        if ($current -eq $null -or $current -eq "") {
            $cell.Value = "nan"
        }
    }
}
